{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\\u00A9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\n\n// Find the index of the first target paragraph so we also remove the\n// blank paragraph that immediately precedes the pair (the paragraph that\n// used to separate the page-break paragraph above from this block).\nlet firstTargetIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (targetTexts.indexOf(items[i].text) !== -1) {\n    firstTargetIndex = i;\n    break;\n  }\n}\n\nconst toDelete = [];\nif (firstTargetIndex > 0 && items[firstTargetIndex - 1].text === \"\") {\n  toDelete.push(items[firstTargetIndex - 1]);\n}\nfor (let i = 0; i < items.length; i++) {\n  if (targetTexts.indexOf(items[i].text) !== -1) {\n    toDelete.push(items[i]);\n  }\n}\n\nfor (const para of toDelete) {\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$targets = @(\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    ([char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\")\n)\n\n# Locate the paragraphs that hold the target text.\n$count = $d.Paragraphs.Count\n$indicesToDelete = New-Object System.Collections.ArrayList\n\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $clean = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($targets -contains $clean) {\n        [void]$indicesToDelete.Add($i)\n    }\n}\n\n# Also remove the blank paragraph immediately preceding the first match\n# (it used to separate the page-break paragraph above from this block).\nif ($indicesToDelete.Count -gt 0) {\n    $first = $indicesToDelete[0]\n    if ($first -gt 1) {\n        $prev = $d.Paragraphs.Item($first - 1)\n        $prevClean = $prev.Range.Text.TrimEnd([char]13, [char]7)\n        if ($prevClean -eq \"\") {\n            [void]$indicesToDelete.Add($first - 1)\n        }\n    }\n}\n\n# Delete from the highest index down so earlier indices stay valid.\n$sorted = $indicesToDelete | Sort-Object -Descending\nforeach ($idx in $sorted) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
